$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.673.71"
$ws.Range("E2").Value = "  +2.80%  "
$ws.Range("D3").Value = "3.206.27"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'599.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.52%  "
$ws.Range("D6").Value = "'155.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.47%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.202.31"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").Value = "'0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").Value = "'0.519"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.54%  "
$ws.Range("E13").Value = "  +3.77%  "
$ws.Range("D14").Value = "'39.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.19%  "
$ws.Range("D15").Value = "3.734.93"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "66.634.53"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "'7.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.31%  "
$ws.Range("D18").Value = "3.211.71"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "'515.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").Value = "'15.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.21%  "
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("D23").Value = "'8.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.42%  "
$ws.Range("D24").Value = "'15.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'85.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'9.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.38%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'3.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.19%  "
$ws.Range("E29").Value = "  +6.07%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +16.65%  "
$ws.Range("E31").Value = "  +4.63%  "
$ws.Range("D32").Value = "'28.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("E33").Value = "  +3.45%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'6.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "'497.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.19%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'54.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("E41").Value = "  +5.96%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0682"
$ws.Range("E42").Value = "  +18.45%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.302"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.39%  "
$ws.Range("D45").Value = "2.935.32"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "'2.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").Value = "'28.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +4.59%  "
$ws.Range("E51").Value = "  +10.00%  "
